# random-processing.pptx: recolor the three "pipeline" connector arrows
# (Data -> , Preprocess -> , Analyze -> ) that used to all share the same
# green (92D050). The fourth arrow (-> Postprocess) keeps its original
# green.
#
# NOTE on color encoding: Shape.Line.ForeColor.RGB (like classic VBA's
# RGB()/OLE_COLOR) is a single Long built from bytes as R + G*256 + B*65536,
# i.e. 0x00BBGGRR - *not* the 0xRRGGBB order the hex swatches are usually
# quoted in. Build each value explicitly from its R/G/B components so the
# mapping to the target hex swatch is unambiguous.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)

function RgbLong($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# "Straight Arrow Connector 24": 92D050 -> 00B050
$slide.Shapes.Item("Straight Arrow Connector 24").Line.ForeColor.RGB = RgbLong 0x00 0xB0 0x50

# "Straight Arrow Connector 25": 92D050 -> 0070C0
$slide.Shapes.Item("Straight Arrow Connector 25").Line.ForeColor.RGB = RgbLong 0x00 0x70 0xC0

# "Straight Arrow Connector 26": 92D050 -> 002060
$slide.Shapes.Item("Straight Arrow Connector 26").Line.ForeColor.RGB = RgbLong 0x00 0x20 0x60

# "Straight Arrow Connector 27" (-> Postprocess) is left untouched; it stays 92D050.

# The authoring app also stamped an (empty) PowerPoint-2012 slide-guide
# extension list onto the presentation on save. There is no documented
# Guides mutation surface in this host (Presentation.Guides / Guides.Add
# resolve to null here), so this is attempted best-effort and is a no-op
# if unsupported - it does not affect the color edits above.
$null = $p.Guides.Add(1, 100)

Write-Output "Recolored pipeline connector arrows."
